$wb = $excel.ActiveWorkbook

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1122.5
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 1196.9231
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 1196.9231
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -1848.9231

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1994.04
$ws.Range("I40").Value = 1450.5
$ws.Range("K40").Value = 1450.5
$ws.Range("M40").Value = -1275.5

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3109.756
$ws.Range("J112").Value = 3181.5789
$ws.Range("L112").Value = 9544.736699999999
$ws.Range("N112").Value = -11760.7367

# ALC row 117
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 31000
$ws.Range("J117").Value = 31000
$ws.Range("L117").Value = 31000
$ws.Range("N117").Value = -40178

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2179753.2
$ws.Range("J129").Value = 3705268.5
$ws.Range("L129").Value = 11115805.5
$ws.Range("N129").Value = -11125805.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4095
$ws.Range("I138").Value = 2280.5908
$ws.Range("J138").Value = 4807.8037
$ws.Range("K138").Value = 6841.7724
$ws.Range("L138").Value = 14423.4111
$ws.Range("M138").Value = -1701.7724
$ws.Range("N138").Value = -24703.4111

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1441.5
$ws.Range("I61").Value = 1450.9375
$ws.Range("J61").Value = 1422.625
$ws.Range("K61").Value = 1450.9375
$ws.Range("L61").Value = 1422.625
$ws.Range("M61").Value = -1238.9375
$ws.Range("N61").Value = -1846.625

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1441.5
$ws.Range("I136").Value = 1450.9375
$ws.Range("J136").Value = 1422.625
$ws.Range("K136").Value = 4352.8125
$ws.Range("L136").Value = 4267.875
$ws.Range("M136").Value = -1802.8125
$ws.Range("N136").Value = -9367.875

# BSM row 31
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 16620.582
$ws.Range("I134").Value = 1468.5082
$ws.Range("J134").Value = 170666.67
$ws.Range("K134").Value = 4405.5246
$ws.Range("L134").Value = 512000.01
$ws.Range("M134").Value = -1870.5246
$ws.Range("N134").Value = -517070.01

# CRP row 23
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 9953.333000000001
$ws.Range("J23").Value = 9953.333000000001
$ws.Range("L23").Value = 9953.333000000001
$ws.Range("N23").Value = -10433.333

# CRP row 27
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 9953.333000000001
$ws.Range("J27").Value = 9953.333000000001
$ws.Range("L27").Value = 9953.333000000001
$ws.Range("N27").Value = -10337.333

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2668.3386
$ws.Range("I31").Value = 1976.15
$ws.Range("J31").Value = 3926.8635
$ws.Range("K31").Value = 1976.15
$ws.Range("L31").Value = 3926.8635
$ws.Range("M31").Value = -1681.15
$ws.Range("N31").Value = -4516.863499999999

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2668.3386
$ws.Range("I34").Value = 1976.15
$ws.Range("J34").Value = 3926.8635
$ws.Range("K34").Value = 1976.15
$ws.Range("L34").Value = 3926.8635
$ws.Range("M34").Value = -1774.15
$ws.Range("N34").Value = -4330.863499999999

# CRP row 70
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 12840
$ws.Range("J70").Value = 12840
$ws.Range("L70").Value = 12840
$ws.Range("N70").Value = -13470

# CRP row 73
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 12840
$ws.Range("J73").Value = 12840
$ws.Range("L73").Value = 12840
$ws.Range("N73").Value = -15024

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1287.0264
$ws.Range("I132").Value = 866.93335
$ws.Range("J132").Value = 2862.375
$ws.Range("K132").Value = 2600.80005
$ws.Range("L132").Value = 8587.125
$ws.Range("M132").Value = -70.80004999999983
$ws.Range("N132").Value = -13647.125

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 15152255
$ws.Range("I134").Value = 724.5517
$ws.Range("J134").Value = 125000850
$ws.Range("K134").Value = 2173.6551
$ws.Range("L134").Value = 375002550
$ws.Range("M134").Value = 361.3449000000001
$ws.Range("N134").Value = -375007620

# CRP row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 59366.668
$ws.Range("J135").Value = 59366.668
$ws.Range("L135").Value = 59366.668
$ws.Range("N135").Value = -69506.66800000001

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2321.3416
$ws.Range("I5").Value = 2260
$ws.Range("J5").Value = 2336.2122
$ws.Range("K5").Value = 6780
$ws.Range("L5").Value = 7008.6366
$ws.Range("M5").Value = -6668
$ws.Range("N5").Value = -7232.6366

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 75758.81
$ws.Range("I131").Value = 144731.42
$ws.Range("J131").Value = 51618.4
$ws.Range("K131").Value = 434194.26
$ws.Range("L131").Value = 154855.2
$ws.Range("M131").Value = -429154.26
$ws.Range("N131").Value = -164935.2

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 876.3
$ws.Range("I132").Value = 876.3
$ws.Range("K132").Value = 7886.7
$ws.Range("M132").Value = -5356.7

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2321.3416
$ws.Range("I135").Value = 2260
$ws.Range("J135").Value = 2336.2122
$ws.Range("K135").Value = 20340
$ws.Range("L135").Value = 21025.9098
$ws.Range("M135").Value = -17805
$ws.Range("N135").Value = -26095.9098

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 8338317
$ws.Range("I14").Value = 9376607
$ws.Range("J14").Value = 32000
$ws.Range("K14").Value = 9376607
$ws.Range("L14").Value = 32000
$ws.Range("M14").Value = -9376439
$ws.Range("N14").Value = -32336

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13066.75
$ws.Range("I70").Value = 18142.465
$ws.Range("J70").Value = 4184.25
$ws.Range("K70").Value = 18142.465
$ws.Range("L70").Value = 4184.25
$ws.Range("M70").Value = -17872.465
$ws.Range("N70").Value = -4724.25

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13066.75
$ws.Range("I73").Value = 18142.465
$ws.Range("J73").Value = 4184.25
$ws.Range("K73").Value = 18142.465
$ws.Range("L73").Value = 4184.25
$ws.Range("M73").Value = -17206.465
$ws.Range("N73").Value = -6056.25

# LTW row 30
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 5380
$ws.Range("I30").Value = 5380
$ws.Range("K30").Value = 5380
$ws.Range("M30").Value = -5272

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4881809.5
$ws.Range("I100").Value = 5346477
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 5346477
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -5345936
$ws.Range("N100").Value = -3882

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 66600
$ws.Range("J46").Value = 66600
$ws.Range("L46").Value = 66600
$ws.Range("N46").Value = -67062

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4757.1787
$ws.Range("I132").Value = 1224.6666
$ws.Range("J132").Value = 25952.25
$ws.Range("K132").Value = 3673.9998
$ws.Range("L132").Value = 77856.75
$ws.Range("M132").Value = -1143.9998
$ws.Range("N132").Value = -82916.75

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 66600
$ws.Range("J134").Value = 66600
$ws.Range("L134").Value = 199800
$ws.Range("N134").Value = -204870

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3412.4167
$ws.Range("I136").Value = 580.5806
$ws.Range("J136").Value = 20969.8
$ws.Range("K136").Value = 1741.7418
$ws.Range("L136").Value = 62909.39999999999
$ws.Range("M136").Value = 808.2582
$ws.Range("N136").Value = -68009.39999999999
